$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Review date" column (B) stores plain text dates (e.g. 2011-03-10).
# Re-generating the list shifted every date whose day-of-month is the 10th
# to the 11th (e.g. 2011-03-10 -> 2011-03-11), leaving all other dates
# (different day-of-month) untouched.

$used = $ws.UsedRange
$rowCount = $used.Rows.Count

for ($r = 2; $r -le $rowCount; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $old = $cell.Text

    if ($old -match '^(\d{4}-\d{2})-10$') {
        $new = $Matches[1] + "-11"

        # Writing a date-shaped string straight into a General-formatted
        # cell would make Excel auto-convert it into a date serial value.
        # Temporarily mark the cell as Text so the literal string is kept,
        # then drop the formatting override again so the cell ends up
        # exactly as it started (General, no explicit style).
        $cell.NumberFormat = "@"
        $cell.Value = $new
        $cell.ClearFormats()
    }
}
